# Adds "Samples" and "Files" tabs/rows to the startup worksheet, and updates
# the "StatQuery" column text (used by all three tabs) to the new combined
# count query.
#
# NOTE: the order in which NEW (not-yet-seen) string values are first
# assigned to cells controls their position in the shared-strings table, so
# cell values below are deliberately written in a specific sequence:
#   9  SamplesTab
#   10 FilesTab
#   11 Samples query
#   12 Files query
#   13 Count query (new text, re-used by rows 2, 3 and 4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New shared query / label text -----------------------------------------

$samplesQuery = @"
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN  ["Standard AC (3week cycles)"]
WITH
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS ``Sample ID``,
            ss.study_subject_id AS ``Case ID``,
            p.program_acronym AS ``Program Code``,
            s.study_acronym AS ``Arm``,
            ss.disease_subtype AS ``Diagnosis``,
            samp.tissue_type AS ``Tissue Type``,
            samp.composition AS ``Tissue Composition``,
            samp.sample_anatomic_site AS ``Sample Anatomic Site``,
            samp.method_of_sample_procurement AS ``Sample Procurement Method``
"@.TrimEnd("`r","`n")

$filesQuery = @"
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE tp.chemotherapy_regimen IN  ["Standard AC (3week cycles)"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS ``File Name``,
    head(labels(samp)) AS ``Association``,
    f.file_description AS ``Description``,
    f.file_format AS ``File Format``,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS ``Program Code``,
    s.study_acronym AS ``Arm``,
    ss.study_subject_id AS ``Case ID``,
    samp.sample_id AS ``Sample ID``
    order by f.file_name
"@.TrimEnd("`r","`n")

$countQuery = @"
MATCH (ss:study_subject)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen in  ["Standard AC (3week cycles)"]
WITH ss
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (ss)<-[:sample_of_study_subject]-(samp)
MATCH (samp)<-[:file_of_sample]-(f)
MATCH (lp)<-[:file_of_laboratory_procedure]-(f)
RETURN COUNT(DISTINCT p) AS Programs,
COUNT(DISTINCT s) AS Arms,
COUNT(DISTINCT ss) AS Cases,
COUNT(DISTINCT samp) AS Samples,
COUNT(DISTINCT lp) AS Assays,
COUNT(DISTINCT f) AS Files
"@.TrimEnd("`r","`n")

$neo4jFile = "TC11_Bento_Filter_Chemo-StandardAC3Wk_Neo4jData.xlsx"
$webFile   = "TC11_Bento_Filter_Chemo-StandardAC3Wk_WebData.xlsx"

# --- Add "Samples" row (row 3) labels/queries first, so the shared-string --
# --- table gets the labels/queries in the order they appear in the file  --

$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B3").Value = $samplesQuery
$ws.Range("B4").Value = $filesQuery

# --- Update existing "Cases" row (row 2) StatQuery column ------------------
# (this introduces the new / final shared string, #13)

$ws.Range("C2").Value = $countQuery

# --- Finish populating rows 3 and 4 -----------------------------------------

$ws.Range("C3").Value = $countQuery
$ws.Range("D3").Value = $neo4jFile
$ws.Range("E3").Value = $webFile

$ws.Range("C4").Value = $countQuery
$ws.Range("D4").Value = $neo4jFile
$ws.Range("E4").Value = $webFile

# --- Wrap text on the query columns, matching the "Normal 2" style ---------

$ws.Range("B2:C4").WrapText = $true

# --- Row heights -------------------------------------------------------------

$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 345.6
$ws.Rows.Item(4).RowHeight = 409.6

# --- Selection / view tweaks to mirror the final saved state ---------------

$ws.Range("D4").Select()
